# STUDY - 4 ThreadLocal学习
# Add a new "ThreadLocal" worksheet after "RabbitMQ" and fill in the new
# study-log rows, then tidy up the view/selection state on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the new worksheet right after the existing "RabbitMQ" sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ThreadLocal"

# Fill in the data. The order in which new string values are first
# assigned controls the shared-string table order, so write the cells in
# this specific sequence (dates first where convenient, then labels).
$ws2.Range("A1").Value = 44043
$ws2.Range("B2").Value = "学习ThreadLocal,使用Apache2"
$ws2.Range("C2").Value = "花费6个小时左右"
$ws2.Range("A2").Value = 44044
$ws2.Range("E2").Value = "windows 中安装Apache2，并使用"
$ws2.Range("D2").Value = "周六"
$ws2.Range("A3").Value = 44045
$ws2.Range("C3").Value = "花费6小时"
$ws2.Range("D1").Value = "周五"
$ws2.Range("B1").Value = "学习简介"
$ws2.Range("B3").Value = "学习ThreadLocal源码，子编写ThreadLocal"
$ws2.Range("D3").Value = "周日"

# Match the date formatting used on the RabbitMQ sheet (column A) by
# copying the existing date cell's format onto the new date cells,
# rather than assigning a NumberFormat string (which would create a
# brand-new custom number format instead of reusing the built-in one).
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore explicit values (PasteSpecial(Formats) should not have touched
# them, but make sure nothing was clobbered).
$ws2.Range("A1").Value = 44043
$ws2.Range("A2").Value = 44044
$ws2.Range("A3").Value = 44045

# On the old sheet, the saved selection now spans the full used rows
# (A1:XFD4) instead of a single cell, and it is no longer the active tab.
$ws1.Rows("1:4").Select() | Out-Null

# Select B1 on the new sheet to match the saved cursor position, then
# make it the active sheet/tab (must happen last so it stays active).
$ws2.Range("B1").Select() | Out-Null
$ws2.Activate() | Out-Null
